$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.794.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.037"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +2.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.46"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.036"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4419"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3799"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07461"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8854"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.80"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.882.19"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -11.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.561"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.758"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07231"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.76"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.32%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009136"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.035"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.829.99"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.322"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.43"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.82%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.64"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.89"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.996"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.333"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.61"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09110"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7779"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.218"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.052"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +8.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.591"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.165"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.61%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05364"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5214"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.841"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.94%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.907"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.729"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.82"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.65"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.726"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4714"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06438"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.888"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.02"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.61"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.15%  "
